$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix goodness scoring for the "Depression" row
$ws.Range("B5").Value = "<1"
$ws.Range("C5").Value = ">0"
$ws.Range("D5").Value = -1
$ws.Range("E5").Value = ">0"

$wb.Save()
